$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# These cells hold text-typed numeric-looking figures (shared strings, not
# real numbers), so the leading apostrophe keeps them stored as text
# (quote-prefixed) instead of being coerced into numeric cells.

# Enterprises density (per 1000 people)
$ws.Range("B13").Value = "'12.86"
$ws.Range("C13").Value = "'0.65"
$ws.Range("D13").Value = "'13.52"

# Employment (% of total)
$ws.Range("B14").Value = "'67.96"
$ws.Range("C14").Value = "'18.66"
$ws.Range("D14").Value = "'86.61"

# Enterprises (% of total)
$ws.Range("B16").Value = "'94.45"
$ws.Range("C16").Value = "'4.79"
$ws.Range("D16").Value = "'99.24"
